$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("index")

# Row 56 ("6_3_rmarkdown" slide deck): the Rmarkdown slides are no longer
# "pending" and now also have exercises.
$ws.Cells.Item(56, 4).Value = "Rmarkdown"
$ws.Cells.Item(56, 6).Value = "yes"

# New row 57: a follow-up "More Rmarkdown" session/slide deck + exercises.
$ws.Cells.Item(57, 1).Value = 16
$ws.Cells.Item(57, 2).Value = 6
$ws.Cells.Item(57, 3).Value = "Miscellanea"
$ws.Cells.Item(57, 4).Value = "More Rmarkdown"
$ws.Cells.Item(57, 5).Value = "6_4_more_rmarkdown"
$ws.Cells.Item(57, 6).Value = "yes"

# Bring the newly added row into view, matching the author's final selection.
$ws.Activate()
$ws.Range("D57").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
